# "Generate Report for Archive"
#
# Changes applied:
#  1. Status value "Ready for handoff" -> "In Translation" on every sheet
#     that carries it (Overview!E2:F2, zh-cn!C2, de-de!C2).
#  2. Narrow the "Status" column(s) from ~17.22 chars to ~13.41 chars:
#       - Overview sheet: columns E and F
#       - zh-cn sheet:    column C
#       - de-de sheet:    column C

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# 1. Update the status text wherever it appears.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# 2. Narrow the status columns (~17.22 -> ~13.41 character width).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
